# Update the DTCP results sheet: rows 2-11 get new Application No / District /
# Approval Type / Permit Issue Date / Total Fees values, and the free-text /
# contact columns (Project Title, Applicant-Owner Signature, Registered
# Engineer Name-Address, Registered Engineer Mail, Registered Engineer Phone)
# are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = "0DJC0CIY"; B = "Ariyalur"; C = "Layout Approval"; D = "22/01/2025"; E = "26,800.00" },
    @{ Row = 3;  A = "ETVG127A"; B = "Ariyalur"; C = "Layout Approval"; D = "22/01/2025"; E = "56,864.00" },
    @{ Row = 4;  A = "YSF78WZR"; B = "Ariyalur"; C = "Layout Approval"; D = "22/01/2025"; E = "62,463.00" },
    @{ Row = 5;  A = "PIEKDXZF"; B = "Ariyalur"; C = "Layout Approval"; D = "02/04/2025"; E = "61,475.50" },
    @{ Row = 6;  A = "99A0FQBY"; B = "Ariyalur"; C = "Layout Approval"; D = "29/04/2025"; E = "46,192.00" },
    @{ Row = 7;  A = "R0B373HM"; B = "Ariyalur"; C = "Layout Approval"; D = "29/04/2025"; E = "28,450.00" },
    @{ Row = 8;  A = "6P4NV2FX"; B = "Ariyalur"; C = "Layout Approval"; D = "09/05/2025"; E = "13,200.00" },
    @{ Row = 9;  A = "9JFTJMKN"; B = "Ariyalur"; C = "Layout Approval"; D = "28/05/2025"; E = "68,851.00" },
    @{ Row = 10; A = "JFCKGKOO"; B = "Ariyalur"; C = "Layout Approval"; D = "28/05/2025"; E = "19,225.50" },
    @{ Row = 11; A = "YRCZUMYR"; B = "Ariyalur"; C = "Layout Approval"; D = "10/07/2025"; E = "16,950.00" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    # Leading apostrophe forces these to stay plain text instead of being
    # auto-parsed into a date serial / number by Excel's smart entry logic
    # (the source data keeps them as literal DD/MM/YYYY and comma-grouped
    # strings).
    $ws.Cells.Item($r, 4).Value = "'" + $item.D
    $ws.Cells.Item($r, 5).Value = "'" + $item.E

    # Clear Project Title / Applicant-Owner Signature / Registered Engineer
    # Name-Address / Registered Engineer Mail / Registered Engineer Phone.
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = ""
    $ws.Cells.Item($r, 9).Value = ""
    $ws.Cells.Item($r, 10).Value = ""
}
